# Actividad 01 - Identificación Hallazgos BD
# Register weaknesses #5, #6 and #7 (rows 6-8) on the "Debilidades" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Debilidades")
$ws.Activate()

# --- Weakness 5 (row 6): "producto.id" missing primary key ---
$ws.Range("B6").Value = "En la base de datos ""ISW1012"" en el esquema ""actividad01"" en la tabla ""producto"" no se identifica la definición de restricción de clave primaria (primary key) sobre la columna ""id""."
$ws.Range("C6").Value = "Integridad de Entidad - Llave Primaria"

# --- Weakness 6 (row 7): "producto.codigo" missing uniqueness constraint ---
$ws.Range("B7").Value = "En la base de datos ""ISW1012"", dentro del esquema ""actividad01"", en la tabla ""producto"" no se identifica la definición de una restricción de unicidad sobre la columna ""codigo"", esto permite la existencia de productos con códigos duplicados."
$ws.Range("C7").Value = "Integridad de Entidad - Llave Única"

# --- Weakness 7 (row 8): "producto.nombre" allows null/blank values ---
$ws.Range("C8").Value = "Integridad de Atributo"
$ws.Range("B8").Value = "En la base de datos ""ISW1012"", dentro del esquema ""actividad01"", en la tabla ""producto"" específicamente en la columna ""nombre"" se permite el ingreso de valores nulos y en blanco, permitiendo la perdida de un dato esencial para la identificación y descripción de los productos."

# Resize the newly-filled rows to fit the wrapped text (matches Excel's
# auto-fit for the column B / C widths used in this sheet).
$ws.Rows.Item(6).RowHeight = 43.2
$ws.Rows.Item(7).RowHeight = 57.6
$ws.Rows.Item(8).RowHeight = 57.6

# The weakness column (B) is vertically centered across the whole table
# (rows 4-19 previously sat on top-alignment; bringing them in line with
# rows 2-3 and the newly filled rows).
$ws.Range("B4:B19").VerticalAlignment = -4108

# Drop the stray formatting-only column F (no longer used).
$ws.Columns.Item(6).Delete()

# Restore the view: scrolled down a bit with D8 as the active cell.
$ws.Range("D8").Select()
